$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.744.70"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.81%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "3.440.87"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.50%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'575.60"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.14%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'160.54"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.30%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "3.444.07"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.53%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "  +9.09%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'7.35"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.00%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "  +2.19%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.440"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.75%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "4.037.70"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.60%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = "  -2.11%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value = "  +4.47%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'28.35"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.67%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "64.762.51"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.72%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "3.490.10"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.97%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'6.38"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.87%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'14.27"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.15%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'386.61"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.27%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'8.19"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.53%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'73.34"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.92%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.546"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.44%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "  +0.12%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "  +14.97%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'9.78"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.67%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'0.181"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.27%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.04%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "  +6.68%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "  +3.63%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "  -0.01%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'6.58"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.76%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'23.65"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.35%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'0.999"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.14%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'7.12"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.02%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("B37").Value = "Monero"
$ws.Range("B37").Style = "Normal"
$ws.Range("C37").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("C37").Style = "Normal"
$ws.Range("D37").Value = "'163.44"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.23%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("B38").Value = "ImmutableX"
$ws.Range("B38").Style = "Normal"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("C38").Style = "Normal"
$ws.Range("D38").Value = "'1.51"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.81%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "3.018.04"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.27%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "  +1.76%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.0767"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.01%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'27.28"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.10%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'4.56"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.10%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("B44").Value = "VeChain"
$ws.Range("B44").Style = "Normal"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("C44").Style = "Normal"
$ws.Range("D44").Value = "'0.0318"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.47%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("B45").Value = "OKB"
$ws.Range("B45").Style = "Normal"
$ws.Range("C45").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("C45").Style = "Normal"
$ws.Range("D45").Value = "'42.83"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.24%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.774"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.74%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'24.66"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +9.29%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "  -0.37%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.882"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +6.78%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "  +4.13%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "  +3.68%  "
$ws.Range("E51").Style = "Normal"
